# Mitarbeiter - Minijobber.xlsx
# "Alter bei AN-Anteil PV hinzugefügt und notwendige Änderungen vorgenommen."
#
# A new question row ("juenger als 23 oder vor 1940 geboren?" / "nein") is
# inserted right above the existing "wohnhaft Sachsen?" row (old row 41),
# which pushes every row from the old row 41 onward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at position 41 - this shifts the old rows 41..48 down to
# 42..49 and carries their formatting/styles/validations along for the ride,
# matching the diff's row renumbering exactly.
$ws.Rows.Item(41).EntireRow.Insert()

# Populate the freshly inserted row with the new question/answer pair.
$ws.Range("A41").Value = "juenger als 23 oder vor 1940 geboren?"
$ws.Range("B41").Value = "nein"

# Restore the view/selection state recorded in the saved workbook.
$ws.Range("A24").Select()
